$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.326.03'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.354.71'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.12%  '
$ws.Range('E3').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '190.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.07%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '557.99'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.349.68'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.24%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.180'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.64'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('E12').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.883.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.57'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '594.50'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.08%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.338.83'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('E17').ClearFormats()

$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('B18').ClearFormats()
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C18').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.359.87'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.21%  '
$ws.Range('E18').ClearFormats()

$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('B19').ClearFormats()
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C19').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.94'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('E19').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.05'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.900'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.46%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.03'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.36%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E26').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.50'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.54'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.72'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.72'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +6.90%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.85'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.61%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '585.10'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.12%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.95'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.104'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.760.30'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.46%  '
$ws.Range('E37').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.68'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.44'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.39%  '
$ws.Range('E40').ClearFormats()

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'PEPE'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0701'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('E41').ClearFormats()

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.127'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.18'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.58%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.38%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.340'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('E46').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -12.84%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0417'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('E48').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.56'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.11%  '
$ws.Range('E51').ClearFormats()
